# Research.xlsx – minimal text change (see commit: "algumas alteracoes de
# textos minimas").
#
# The shared-strings table lost the string "Pedra de Fogo" and gained a new
# string "Pedra da Cidade" appended at the end. Since every other cell in
# column C (rows 4-11) still shows the exact same Portuguese text it showed
# before, the only real content edit is the translated name in cell C3:
#   "Pedra de Fogo"  ->  "Pedra da Cidade"
# (the cascading shared-string index shuffle seen in the raw xlsx diff is
# just a side effect of removing/appending to the shared strings table).
#
# The sheet view also ended up with the cursor resting on C8 (instead of
# the former C3:C11 selection), and column D was widened/auto-fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The actual content edit.
$ws.Range("C3").Value = "Pedra da Cidade"

# Column D ended up auto-fit to its (unchanged) contents in the authored
# file; reproduce that sizing action here.
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null

# Final cursor position recorded in the sheet view.
$ws.Range("C8").Select()
